# Reorders the weekly data rows (rows 2-19) of the "Hortaliza, Terminal La
# Palmera de La Serena - Alcachofa" sheet. Only the row each record lives in
# changes; the fields/values of each record stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # column A
$lastCol  = 18  # column R

# Mapping of target row -> source row (1-based worksheet rows).
$rowMap = @{
    2  = 7
    3  = 9
    4  = 18
    5  = 2
    6  = 10
    7  = 11
    8  = 15
    9  = 3
    10 = 5
    11 = 13
    12 = 14
    13 = 8
    14 = 6
    15 = 19
    16 = 4
    17 = 12
    18 = 16
    19 = 17
}

# 1) Snapshot every source row's values before any writes happen, so that
#    overwriting one row doesn't clobber data still needed for another.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowValues = @{}
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $rowValues[$col] = $ws.Cells.Item($srcRow, $col).Value()
        }
        $snapshot[$srcRow] = $rowValues
    }
}

# 2) Write the snapshotted values into their new (target) row positions.
foreach ($targetRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$targetRow]
    $rowValues = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $rowValues[$col]
    }
}
